$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 26
$endRow = 33

$A = @("yield","yield","yield","yield","yield","yield","yield","yield")
$B = @(1,2,2,3,3,4,4,4)
$C = @("D03_cuero","D09_T_liquor","D09_Tri_FrecA","D10_T_out_gelatin","D10_T_out_water","D12_P_atom_regn","D12_P_filt","D12_P_regen")
$D = @(1,2,3,4,5,6,7,8)
$E = @(
    "Leather as raw material",
    "Liqueur temperature paraflash (ºC)",
    "Triplex pump A frecuency (osc/min)",
    "Sterilizer gelatin outflow temperature (ºC)",
    "Sterilizer water outflow temperature (ºC)",
    "Regenerator atomization pressure (psi)",
    "Filter drop pressure (in H2O)",
    "Regenerator drop pressure (in H2O)"
)
$F = @(1,50.277053140096598,20.8319082125603,48.176158204643897,40.110690315408696,0.19253048079134999,0.29206314699792901,0.87258482861743702)
$G = @("text","number","number","number","number","number","number","number")
$H = @("input_text","input_text","input_text","input_text","input_text","input_text","input_text","input_text")
$I = $F
$J = @("True","True","True","True","True","True","True","True")

for ($r = $startRow; $r -le $endRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $A[$r - $startRow]
}
for ($r = $startRow; $r -le $endRow; $r++) {
    $ws.Cells.Item($r, 2).Value = $B[$r - $startRow]
}
for ($r = $startRow; $r -le $endRow; $r++) {
    $ws.Cells.Item($r, 3).Value = $C[$r - $startRow]
}
for ($r = $startRow; $r -le $endRow; $r++) {
    $ws.Cells.Item($r, 4).Value = $D[$r - $startRow]
}
for ($r = $startRow; $r -le $endRow; $r++) {
    $ws.Cells.Item($r, 5).Value = $E[$r - $startRow]
}
for ($r = $startRow; $r -le $endRow; $r++) {
    $ws.Cells.Item($r, 6).Value = $F[$r - $startRow]
}
for ($r = $startRow; $r -le $endRow; $r++) {
    $ws.Cells.Item($r, 7).Value = $G[$r - $startRow]
}
for ($r = $startRow; $r -le $endRow; $r++) {
    $ws.Cells.Item($r, 8).Value = $H[$r - $startRow]
}
for ($r = $startRow; $r -le $endRow; $r++) {
    $ws.Cells.Item($r, 9).Value = $I[$r - $startRow]
}
for ($r = $startRow; $r -le $endRow; $r++) {
    $ws.Cells.Item($r, 10).Value = $J[$r - $startRow]
}

$ws.Range("C26:C33,E26:E33").Style = $ws.Range("C17").Style
$ws.Range("D26:D33,F26:F33").Style = $ws.Range("D17").Style

$ws.Application.ActiveWindow.ScrollRow = 12
$ws.Range("B26").Select()
